$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D retains exact text formatting (avoid numeric auto-conversion)
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "20.565.16"
$ws.Range("E2").Value = "  +1.77%  "
$ws.Range("D3").Value = "1.474.35"
$ws.Range("E3").Value = "  +2.51%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "0.9588"
$ws.Range("E5").Value = "  +3.90%  "
$ws.Range("D6").Value = "277.18"
$ws.Range("E6").Value = "  +1.24%  "
$ws.Range("D7").Value = "0.3521"
$ws.Range("E7").Value = "  -3.16%  "
$ws.Range("D8").Value = "0.3070"
$ws.Range("E8").Value = "  +0.56%  "
$ws.Range("D9").Value = "1.086"
$ws.Range("E9").Value = "  +7.18%  "
$ws.Range("D10").Value = "39.42"
$ws.Range("E10").Value = "  +0.87%  "
$ws.Range("D11").Value = "0.06635"
$ws.Range("E11").Value = "  +2.67%  "
$ws.Range("E12").Value = "  +0.37%  "
$ws.Range("D13").Value = "18.09"
$ws.Range("E13").Value = "  +4.56%  "
$ws.Range("D14").Value = "5.452"
$ws.Range("E14").Value = "  +2.66%  "
$ws.Range("D15").Value = "6.170"
$ws.Range("E15").Value = "  +2.54%  "
$ws.Range("D16").Value = "0.9596"
$ws.Range("E16").Value = "  +1.80%  "
$ws.Range("D17").Value = "0.00001018"
$ws.Range("E17").Value = "  +1.18%  "
$ws.Range("D18").Value = "1.472.76"
$ws.Range("E18").Value = "  +2.35%  "
$ws.Range("D19").Value = "0.05964"
$ws.Range("E19").Value = "  +5.58%  "
$ws.Range("D20").Value = "69.03"
$ws.Range("E20").Value = "  +1.82%  "
$ws.Range("D21").Value = "5.485"
$ws.Range("E21").Value = "  +2.73%  "
$ws.Range("E22").Value = "  +2.39%  "
$ws.Range("D23").Value = "11.18"
$ws.Range("E23").Value = "  +4.06%  "
$ws.Range("D24").Value = "2.277"
$ws.Range("E24").Value = "  +1.45%  "
$ws.Range("D25").Value = "20.582.55"
$ws.Range("E25").Value = "  +1.65%  "
$ws.Range("D26").Value = "146.16"
$ws.Range("E26").Value = "  +4.85%  "
$ws.Range("D27").Value = "2.090"
$ws.Range("E27").Value = "  +3.64%  "
$ws.Range("D28").Value = "17.13"
$ws.Range("E28").Value = "  +1.97%  "
$ws.Range("D29").Value = "1.632.10"
$ws.Range("E29").Value = "  +2.52%  "
$ws.Range("D30").Value = "114.44"
$ws.Range("E30").Value = "  +4.24%  "
$ws.Range("D31").Value = "3.948"
$ws.Range("E31").Value = "  -3.03%  "
$ws.Range("D32").Value = "4.941"
$ws.Range("E32").Value = "  +3.83%  "
$ws.Range("D33").Value = "0.07905"
$ws.Range("E33").Value = "  +3.49%  "
$ws.Range("D34").Value = "0.7963"
$ws.Range("E34").Value = "  +3.48%  "
$ws.Range("D35").Value = "1.205"
$ws.Range("E35").Value = "  +8.50%  "
$ws.Range("D36").Value = "1.431"
$ws.Range("E36").Value = "  -1.36%  "
$ws.Range("D37").Value = "0.05669"
$ws.Range("E37").Value = "  +0.36%  "
$ws.Range("D38").Value = "4.708"
$ws.Range("E38").Value = "  +2.28%  "
$ws.Range("D39").Value = "0.9602"
$ws.Range("E39").Value = "  +3.17%  "
$ws.Range("D40").Value = "0.02020"
$ws.Range("E40").Value = "  +2.41%  "
$ws.Range("D41").Value = "10.27"
$ws.Range("E41").Value = "  +1.91%  "
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").Value = "7.420"
$ws.Range("E42").Value = "  +6.82%  "
$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D43").Value = "0.1847"
$ws.Range("E43").Value = "  +1.26%  "
$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").Value = "0.5223"
$ws.Range("E44").Value = "  +1.37%  "
$ws.Range("B45").Value = "PancakeSwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D45").Value = "3.511"
$ws.Range("E45").Value = "  +1.20%  "
$ws.Range("D46").Value = "12.01"
$ws.Range("E46").Value = "  +2.69%  "
$ws.Range("D47").Value = "119.89"
$ws.Range("E47").Value = "  +5.27%  "
$ws.Range("E48").Value = "  +2.50%  "
$ws.Range("D49").Value = "1.808"
$ws.Range("E49").Value = "  +5.52%  "
$ws.Range("D50").Value = "0.06406"
$ws.Range("E50").Value = "  +0.95%  "
$ws.Range("D51").Value = "0.9937"
$ws.Range("E51").Value = "  +0.42%  "
